$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1998.7826
$ws.Range("J17").Value = 1998.7826
$ws.Range("L17").Value = 5996.3478
$ws.Range("N17").Value = -6332.3478
$ws.Range("H40").Value = 5429.6665
$ws.Range("I40").Value = 6398.143
$ws.Range("J40").Value = 2040
$ws.Range("K40").Value = 6398.143
$ws.Range("L40").Value = 2040
$ws.Range("M40").Value = -6223.143
$ws.Range("N40").Value = -2390
$ws.Range("H116").Value = 2863.262
$ws.Range("I116").Value = 2241.5
$ws.Range("J116").Value = 4417.6665
$ws.Range("K116").Value = 2241.5
$ws.Range("L116").Value = 4417.6665
$ws.Range("M116").Value = 1200.5
$ws.Range("N116").Value = -11301.6665
$ws.Range("H135").Value = 14286464
$ws.Range("I135").Value = 805.0645
$ws.Range("J135").Value = 125000320
$ws.Range("K135").Value = 7245.5805
$ws.Range("L135").Value = 1125002880
$ws.Range("M135").Value = -4710.5805
$ws.Range("N135").Value = -1125007950

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H55").Value = 38000
$ws.Range("J55").Value = 38000
$ws.Range("L55").Value = 38000
$ws.Range("N55").Value = -38630
$ws.Range("H61").Value = 2355.0571
$ws.Range("I61").Value = 1386
$ws.Range("J61").Value = 3381.1177
$ws.Range("K61").Value = 1386
$ws.Range("L61").Value = 3381.1177
$ws.Range("M61").Value = -1174
$ws.Range("N61").Value = -3805.1177
$ws.Range("H74").Value = 1761.1555
$ws.Range("I74").Value = 1502.2572
$ws.Range("J74").Value = 2667.3
$ws.Range("K74").Value = 1502.2572
$ws.Range("L74").Value = 2667.3
$ws.Range("M74").Value = -628.2572
$ws.Range("N74").Value = -4415.3
$ws.Range("H77").Value = 1761.1555
$ws.Range("I77").Value = 1502.2572
$ws.Range("J77").Value = 2667.3
$ws.Range("K77").Value = 7511.286
$ws.Range("L77").Value = 13336.5
$ws.Range("M77").Value = -3143.286
$ws.Range("N77").Value = -22072.5
$ws.Range("H122").Value = 1868.0625
$ws.Range("I122").Value = 1792.1428
$ws.Range("J122").Value = 2399.5
$ws.Range("K122").Value = 5376.428400000001
$ws.Range("L122").Value = 7198.5
$ws.Range("M122").Value = -2926.428400000001
$ws.Range("N122").Value = -12098.5
$ws.Range("H136").Value = 2355.0571
$ws.Range("I136").Value = 1386
$ws.Range("J136").Value = 3381.1177
$ws.Range("K136").Value = 4158
$ws.Range("L136").Value = 10143.3531
$ws.Range("M136").Value = -1608
$ws.Range("N136").Value = -15243.3531

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 4949.6
$ws.Range("I134").Value = 5099.857
$ws.Range("K134").Value = 15299.571
$ws.Range("M134").Value = -12764.571

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5959081
$ws.Range("I31").Value = 2283.3
$ws.Range("J31").Value = 7254037
$ws.Range("K31").Value = 2283.3
$ws.Range("L31").Value = 7254037
$ws.Range("M31").Value = -1988.3
$ws.Range("N31").Value = -7254627
$ws.Range("H34").Value = 5959081
$ws.Range("I34").Value = 2283.3
$ws.Range("J34").Value = 7254037
$ws.Range("K34").Value = 2283.3
$ws.Range("L34").Value = 7254037
$ws.Range("M34").Value = -2081.3
$ws.Range("N34").Value = -7254441
$ws.Range("H58").Value = 1781.0731
$ws.Range("I58").Value = 1158.4584
$ws.Range("K58").Value = 1158.4584
$ws.Range("M58").Value = -955.4584
$ws.Range("H136").Value = 1781.0731
$ws.Range("I136").Value = 1158.4584
$ws.Range("K136").Value = 3475.3752
$ws.Range("M136").Value = -925.3751999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 7723.7144
$ws.Range("J5").Value = 1100
$ws.Range("L5").Value = 3300
$ws.Range("N5").Value = -3524
$ws.Range("H120").Value = 3000000
$ws.Range("J120").Value = 0
$ws.Range("L120").Value = 0
$ws.Range("N120").ClearContents()
$ws.Range("H122").Value = 5341.5
$ws.Range("I122").Value = 460.9
$ws.Range("J122").Value = 9408.666999999999
$ws.Range("K122").Value = 4148.099999999999
$ws.Range("L122").Value = 84678.003
$ws.Range("M122").Value = -1698.099999999999
$ws.Range("N122").Value = -89578.003
$ws.Range("H132").Value = 1845.8948
$ws.Range("I132").Value = 1299.5294
$ws.Range("J132").Value = 6490
$ws.Range("K132").Value = 11695.7646
$ws.Range("L132").Value = 58410
$ws.Range("M132").Value = -9165.764599999999
$ws.Range("N132").Value = -63470
$ws.Range("H135").Value = 7723.7144
$ws.Range("J135").Value = 1100
$ws.Range("L135").Value = 9900
$ws.Range("N135").Value = -14970

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1559.7693
$ws.Range("I122").Value = 1553
$ws.Range("J122").Value = 1597
$ws.Range("K122").Value = 4659
$ws.Range("L122").Value = 4791
$ws.Range("M122").Value = -2209
$ws.Range("N122").Value = -9691
$ws.Range("H132").Value = 4639.3076
$ws.Range("I132").Value = 3005.2
$ws.Range("J132").Value = 5660.625
$ws.Range("K132").Value = 9015.599999999999
$ws.Range("L132").Value = 16981.875
$ws.Range("M132").Value = -6485.599999999999
$ws.Range("N132").Value = -22041.875

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 4499.15
$ws.Range("I46").Value = 839.3333
$ws.Range("J46").Value = 6695.04
$ws.Range("K46").Value = 839.3333
$ws.Range("L46").Value = 6695.04
$ws.Range("M46").Value = -651.3333
$ws.Range("N46").Value = -7071.04
$ws.Range("H122").Value = 85715.836
$ws.Range("I122").Value = 113376.664
$ws.Range("J122").Value = 2733.3333
$ws.Range("K122").Value = 340129.992
$ws.Range("L122").Value = 8199.999899999999
$ws.Range("M122").Value = -337679.992
$ws.Range("N122").Value = -13099.9999
$ws.Range("H132").Value = 3823.4412
$ws.Range("I132").Value = 2605.2632
$ws.Range("K132").Value = 7815.7896
$ws.Range("M132").Value = -5285.7896
$ws.Range("H136").Value = 2070.5652
$ws.Range("I136").Value = 1575.4736
$ws.Range("J136").Value = 4422.25
$ws.Range("K136").Value = 4726.4208
$ws.Range("L136").Value = 13266.75
$ws.Range("M136").Value = -2176.4208
$ws.Range("N136").Value = -18366.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H39").Value = 24247.5
$ws.Range("I39").Value = 0
$ws.Range("K39").Value = 0
$ws.Range("M39").ClearContents()
$ws.Range("H113").Value = 1184.8518
$ws.Range("I113").Value = 1219.6666
$ws.Range("J113").Value = 1141.3334
$ws.Range("K113").Value = 3658.9998
$ws.Range("L113").Value = 3424.0002
$ws.Range("M113").Value = -1488.9998
$ws.Range("N113").Value = -7764.0002
$ws.Range("H122").Value = 1929.7646
$ws.Range("I122").Value = 2021.8572
$ws.Range("J122").Value = 1500
$ws.Range("K122").Value = 6065.571599999999
$ws.Range("L122").Value = 4500
$ws.Range("M122").Value = -3615.571599999999
$ws.Range("N122").Value = -9400
$ws.Range("H123").Value = 43076.332
$ws.Range("J123").Value = 43076.332
$ws.Range("L123").Value = 43076.332
$ws.Range("N123").Value = -52876.332
$ws.Range("H136").Value = 1257.7
$ws.Range("I136").Value = 1062.3478
$ws.Range("J136").Value = 1899.5714
$ws.Range("K136").Value = 3187.0434
$ws.Range("L136").Value = 5698.7142
$ws.Range("M136").Value = -637.0434
$ws.Range("N136").Value = -10798.7142
